# Update the dSF (column F) values for the affected rows, per the
# "repull data, push all data, mean calculation" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 6
    6  = 0
    9  = -3
    13 = -1
    14 = 0
    15 = -1
    18 = 5
    19 = 2
    24 = 0
    32 = -2
    37 = 2
    42 = -2
    43 = 1
    45 = 12
    46 = -7
    47 = 1
    49 = 3
    51 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
